$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.246.61'
$ws.Range("E2").Value = '  -3.44%  '
$ws.Range("D3").Value = '3.507.56'
$ws.Range("E3").Value = '  -5.07%  '
$style_4D = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = $style_4D
$ws.Range("E4").Value = '  -0.15%  '
$style_5D = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.80'
$ws.Range("D5").Style = $style_5D
$ws.Range("E5").Value = '  -1.06%  '
$style_6D = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.94'
$ws.Range("D6").Style = $style_6D
$ws.Range("E6").Value = '  -4.37%  '
$style_7D = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.614'
$ws.Range("D7").Style = $style_7D
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("D8").Value = '3.500.30'
$ws.Range("E8").Value = '  -5.08%  '
$ws.Range("E9").Value = '  -0.07%  '
$style_10D = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.189'
$ws.Range("D10").Style = $style_10D
$ws.Range("E10").Value = '  -6.71%  '
$style_11D = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.70'
$ws.Range("D11").Style = $style_11D
$ws.Range("E11").Value = '  +5.17%  '
$style_12D = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.593'
$ws.Range("D12").Style = $style_12D
$ws.Range("E12").Value = '  -3.42%  '
$style_13D = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.62'
$ws.Range("D13").Style = $style_13D
$ws.Range("E13").Value = '  -6.64%  '
$style_14D = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000275'
$ws.Range("D14").Style = $style_14D
$ws.Range("E14").Value = '  -4.23%  '
$style_15D = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '674.54'
$ws.Range("D15").Style = $style_15D
$ws.Range("E15").Value = '  -0.92%  '
$ws.Range("D16").Value = '4.062.95'
$ws.Range("E16").Value = '  -5.43%  '
$style_17D = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '8.67'
$ws.Range("D17").Style = $style_17D
$ws.Range("E17").Value = '  -4.03%  '
$ws.Range("D18").Value = '69.164.56'
$ws.Range("E18").Value = '  -3.71%  '
$ws.Range("D19").Value = '3.497.15'
$ws.Range("E19").Value = '  -5.48%  '
$ws.Range("E20").Value = '  -1.50%  '
$style_21D = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.38'
$ws.Range("D21").Style = $style_21D
$ws.Range("E21").Value = '  -3.97%  '
$style_22D = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.25'
$ws.Range("D22").Style = $style_22D
$ws.Range("E22").Value = '  -3.35%  '
$style_23D = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.899'
$ws.Range("D23").Style = $style_23D
$ws.Range("E23").Value = '  -4.94%  '
$style_24D = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '16.09'
$ws.Range("D24").Style = $style_24D
$ws.Range("E24").Value = '  -9.64%  '
$style_25D = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.55'
$ws.Range("D25").Style = $style_25D
$ws.Range("E25").Value = '  -5.65%  '
$style_26D = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.85'
$ws.Range("D26").Style = $style_26D
$ws.Range("E26").Value = '  -4.71%  '
$style_27D = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").Style = $style_27D
$ws.Range("E27").Value = '  -0.12%  '
$style_28D = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.65'
$ws.Range("D28").Style = $style_28D
$ws.Range("E28").Value = '  -6.68%  '
$style_29D = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.43'
$ws.Range("D29").Style = $style_29D
$ws.Range("E29").Value = '  -8.71%  '
$style_30D = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.87'
$ws.Range("D30").Style = $style_30D
$ws.Range("E30").Value = '  -7.99%  '
$style_31D = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.67'
$ws.Range("D31").Style = $style_31D
$ws.Range("E31").Value = '  -6.96%  '
$style_32D = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.19'
$ws.Range("D32").Style = $style_32D
$ws.Range("E32").Value = '  -8.14%  '
$ws.Range("E33").Value = '  -5.46%  '
$style_34D = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.24'
$ws.Range("D34").Style = $style_34D
$ws.Range("E34").Value = '  -1.61%  '
$style_35D = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '589.82'
$ws.Range("D35").Style = $style_35D
$ws.Range("E35").Value = '  +4.02%  '
$ws.Range("B36").Value = 'dogwifhat'
$ws.Range("C36").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$style_36D = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.59'
$ws.Range("D36").Style = $style_36D
$ws.Range("E36").Value = '  -14.26%  '
$ws.Range("B37").Value = 'Cosmos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$style_37D = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.85'
$ws.Range("D37").Style = $style_37D
$ws.Range("E37").Value = '  -4.02%  '
$style_38D = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.104'
$ws.Range("D38").Style = $style_38D
$ws.Range("E38").Value = '  -5.27%  '
$style_39D = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '56.99'
$ws.Range("D39").Style = $style_39D
$ws.Range("E39").Value = '  -4.14%  '
$style_40D = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("D40").Style = $style_40D
$ws.Range("E40").Value = '  +0.01%  '
$style_41D = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0436'
$ws.Range("D41").Style = $style_41D
$ws.Range("E41").Value = '  -6.09%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$style_42D = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.136'
$ws.Range("D42").Style = $style_42D
$ws.Range("E42").Value = '  -6.02%  '
$ws.Range("B43").Value = 'TheGraph'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$style_43D = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.334'
$ws.Range("D43").Style = $style_43D
$ws.Range("E43").Value = '  -5.23%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '3.406.65'
$ws.Range("E44").Value = '  -9.07%  '
$style_45D = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '33.25'
$ws.Range("D45").Style = $style_45D
$ws.Range("E45").Value = '  -6.80%  '
$ws.Range("D46").Value = '0.0₃0706'
$ws.Range("E46").Value = '  -9.06%  '
$style_47D = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.88'
$ws.Range("D47").Style = $style_47D
$ws.Range("E47").Value = '  -0.95%  '
$style_48D = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.60'
$ws.Range("D48").Style = $style_48D
$ws.Range("E48").Value = '  -7.76%  '
$ws.Range("E49").Value = '  -0.74%  '
$style_50D = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.46'
$ws.Range("D50").Style = $style_50D
$ws.Range("E50").Value = '  -1.90%  '
$style_51D = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.66'
$ws.Range("D51").Style = $style_51D
$ws.Range("E51").Value = '  +14.89%  '
